$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7951765
$ws.Range("I43").Value = 34000
$ws.Range("K43").Value = 34000
$ws.Range("M43").Value = -33931
$ws.Range("H112").Value = 2853
$ws.Range("J112").Value = 2853
$ws.Range("L112").Value = 8559
$ws.Range("N112").Value = -10775
$ws.Range("H137").Value = 1388.8889
$ws.Range("I137").Value = 1266.1666
$ws.Range("K137").Value = 3798.4998
$ws.Range("M137").Value = -1248.4998
$ws.Range("H138").Value = 783.4545000000001
$ws.Range("I138").Value = 751.6667
$ws.Range("J138").Value = 2500
$ws.Range("K138").Value = 2255.0001
$ws.Range("L138").Value = 7500
$ws.Range("M138").Value = 2884.9999
$ws.Range("N138").Value = -17780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12042.667
$ws.Range("I2").Value = 895.1667
$ws.Range("K2").Value = 895.1667
$ws.Range("M2").Value = -782.1667
$ws.Range("H32").Value = 5945.885
$ws.Range("I32").Value = 5504.304
$ws.Range("J32").Value = 9331.333000000001
$ws.Range("K32").Value = 5504.304
$ws.Range("L32").Value = 9331.333000000001
$ws.Range("M32").Value = -5217.304
$ws.Range("N32").Value = -9905.333000000001
$ws.Range("H108").Value = 5000
$ws.Range("J108").Value = 5000
$ws.Range("L108").Value = 5000
$ws.Range("N108").Value = -12680
$ws.Range("H110").Value = 1847.5454
$ws.Range("I110").Value = 1263.75
$ws.Range("J110").Value = 3404.3333
$ws.Range("K110").Value = 1263.75
$ws.Range("L110").Value = 3404.3333
$ws.Range("M110").Value = 781.25
$ws.Range("N110").Value = -7494.3333
$ws.Range("H116").Value = 12042.667
$ws.Range("I116").Value = 895.1667
$ws.Range("K116").Value = 895.1667
$ws.Range("M116").Value = 1398.8333
$ws.Range("H132").Value = 2484.6765
$ws.Range("I132").Value = 1710.409
$ws.Range("K132").Value = 5131.227000000001
$ws.Range("M132").Value = -2601.227000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 12042.667
$ws.Range("I3").Value = 895.1667
$ws.Range("K3").Value = 895.1667
$ws.Range("M3").Value = -781.1667
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H22").Value = 366.66666
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -127
$ws.Range("N22").Value = -846
$ws.Range("H94").Value = 35715732
$ws.Range("I94").Value = 41667856
$ws.Range("J94").Value = 2999
$ws.Range("K94").Value = 41667856
$ws.Range("L94").Value = 2999
$ws.Range("M94").Value = -41667405
$ws.Range("N94").Value = -3901
$ws.Range("H107").Value = 2344
$ws.Range("I107").Value = 1697
$ws.Range("J107").Value = 2451.8333
$ws.Range("K107").Value = 1697
$ws.Range("L107").Value = 2451.8333
$ws.Range("M107").Value = 223
$ws.Range("N107").Value = -6291.8333
$ws.Range("H131").Value = 30890
$ws.Range("J131").Value = 30890
$ws.Range("L131").Value = 30890
$ws.Range("N131").Value = -40970
$ws.Range("H140").Value = 21689.691
$ws.Range("J140").Value = 21689.691
$ws.Range("L140").Value = 21689.691
$ws.Range("N140").Value = -32049.691

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 319499.75
$ws.Range("I3").Value = 92666.336
$ws.Range("K3").Value = 92666.336
$ws.Range("M3").Value = -92553.336
$ws.Range("H6").Value = 1133.5
$ws.Range("I6").Value = 1060.2
$ws.Range("J6").Value = 1500
$ws.Range("K6").Value = 1060.2
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -947.2
$ws.Range("N6").Value = -1726
$ws.Range("H22").Value = 410
$ws.Range("I22").Value = 316.66666
$ws.Range("K22").Value = 316.66666
$ws.Range("M22").Value = 33.33334000000002
$ws.Range("H31").Value = 1375.8334
$ws.Range("I31").Value = 1291
$ws.Range("J31").Value = 1800
$ws.Range("K31").Value = 1291
$ws.Range("L31").Value = 1800
$ws.Range("M31").Value = -996
$ws.Range("N31").Value = -2390
$ws.Range("H34").Value = 1375.8334
$ws.Range("I34").Value = 1291
$ws.Range("J34").Value = 1800
$ws.Range("K34").Value = 1291
$ws.Range("L34").Value = 1800
$ws.Range("M34").Value = -1089
$ws.Range("N34").Value = -2204
$ws.Range("H62").Value = 9526685
$ws.Range("I62").Value = 2855.2222
$ws.Range("J62").Value = 66669668
$ws.Range("K62").Value = 2855.2222
$ws.Range("L62").Value = 66669668
$ws.Range("M62").Value = -2231.2222
$ws.Range("N62").Value = -66670916
$ws.Range("H65").Value = 9526685
$ws.Range("I65").Value = 2855.2222
$ws.Range("J65").Value = 66669668
$ws.Range("K65").Value = 14276.111
$ws.Range("L65").Value = 333348340
$ws.Range("M65").Value = -11156.111
$ws.Range("N65").Value = -333354580
$ws.Range("H86").Value = 4779649
$ws.Range("I86").Value = 11114579
$ws.Range("K86").Value = 11114579
$ws.Range("M86").Value = -11113456
$ws.Range("H89").Value = 4779649
$ws.Range("I89").Value = 11114579
$ws.Range("K89").Value = 55572895
$ws.Range("M89").Value = -55567279
$ws.Range("H107").Value = 672.06665
$ws.Range("I107").Value = 500.5
$ws.Range("J107").Value = 868.1429000000001
$ws.Range("K107").Value = 500.5
$ws.Range("L107").Value = 868.1429000000001
$ws.Range("M107").Value = 1419.5
$ws.Range("N107").Value = -4708.1429
$ws.Range("H132").Value = 7569.909
$ws.Range("I132").Value = 10879.182
$ws.Range("J132").Value = 4260.636
$ws.Range("K132").Value = 32637.546
$ws.Range("L132").Value = 12781.908
$ws.Range("M132").Value = -30107.546
$ws.Range("N132").Value = -17841.908
$ws.Range("H134").Value = 15153509
$ws.Range("I134").Value = 18520500
$ws.Range("K134").Value = 55561500
$ws.Range("M134").Value = -55558965

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1990
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1990
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 5970
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -6308
$ws.Range("H30").Value = 1990
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1990
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 5970
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -6174
$ws.Range("H131").Value = 10205239
$ws.Range("I131").Value = 166666980
$ws.Range("J131").Value = 1212.3586
$ws.Range("K131").Value = 500000940
$ws.Range("L131").Value = 3637.0758
$ws.Range("M131").Value = -499995900
$ws.Range("N131").Value = -13717.0758

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 10000000
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H24").Value = 1669984.4
$ws.Range("I24").Value = 5001953
$ws.Range("J24").Value = 4000
$ws.Range("K24").Value = 5001953
$ws.Range("L24").Value = 4000
$ws.Range("M24").Value = -5001780
$ws.Range("N24").Value = -4346
$ws.Range("H97").Value = 1019.1818
$ws.Range("I97").Value = 944.44446
$ws.Range("J97").Value = 1355.5
$ws.Range("K97").Value = 944.44446
$ws.Range("L97").Value = 1355.5
$ws.Range("M97").Value = -448.44446
$ws.Range("N97").Value = -2347.5
$ws.Range("H126").Value = 1939.0714
$ws.Range("I126").Value = 1813.3636
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 5440.0908
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -2970.0908
$ws.Range("N126").Value = -12140
$ws.Range("H132").Value = 2386.5715
$ws.Range("I132").Value = 2001.381
$ws.Range("J132").Value = 3542.1428
$ws.Range("K132").Value = 6004.143
$ws.Range("L132").Value = 10626.4284
$ws.Range("M132").Value = -3474.143
$ws.Range("N132").Value = -15686.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 2000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 2000
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -2214
$ws.Range("H46").Value = 2733.6667
$ws.Range("J46").Value = 2733.6667
$ws.Range("L46").Value = 2733.6667
$ws.Range("N46").Value = -3109.6667
$ws.Range("H122").Value = 35717356
$ws.Range("J122").Value = 3255
$ws.Range("L122").Value = 9765
$ws.Range("N122").Value = -14665
$ws.Range("H132").Value = 23975.955
$ws.Range("I132").Value = 1434.826
$ws.Range("J132").Value = 47541.684
$ws.Range("K132").Value = 4304.478
$ws.Range("L132").Value = 142625.052
$ws.Range("M132").Value = -1774.478
$ws.Range("N132").Value = -147685.052
$ws.Range("H136").Value = 8553.4
$ws.Range("I136").Value = 26320.75
$ws.Range("K136").Value = 78962.25
$ws.Range("M136").Value = -76412.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 50698
$ws.Range("J119").Value = 50698
$ws.Range("L119").Value = 50698
$ws.Range("N119").Value = -60374
$ws.Range("H132").Value = 3365.65
$ws.Range("I132").Value = 3065.7144
$ws.Range("J132").Value = 4065.5
$ws.Range("K132").Value = 9197.143199999999
$ws.Range("L132").Value = 12196.5
$ws.Range("M132").Value = -6667.143199999999
$ws.Range("N132").Value = -17256.5
$ws.Range("H135").Value = 60238.332
$ws.Range("J135").Value = 60238.332
$ws.Range("L135").Value = 60238.332
$ws.Range("N135").Value = -70378.33199999999
$ws.Range("H136").Value = 1662.5385
$ws.Range("I136").Value = 1016.4286
$ws.Range("J136").Value = 2416.3333
$ws.Range("K136").Value = 3049.2858
$ws.Range("L136").Value = 7248.999899999999
$ws.Range("M136").Value = -499.2857999999997
$ws.Range("N136").Value = -12348.9999

Write-Output "Applied all Kujata_Profits updates"